# ModificarClienteCorporativo.xlsx
# "update entregable 1, 2 y 3"
#
# The recorded transaction's "Fecha" (date/time stamp) cell G2 is updated
# from the old run's timestamp to the latest run's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "3 jul. 2023, 10:59:01"
